# Apply updated numeric values to Sheet1 as described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.513
$ws.Range("C3").Value = -12.032
$ws.Range("D8").Value = -8.532
$ws.Range("D11").Value = -7.457000000000001
$ws.Range("A12").Value = -21.629
$ws.Range("B14").Value = 6.114
$ws.Range("D14").Value = -7.678
$ws.Range("D15").Value = -8.175000000000001
$ws.Range("B26").Value = 6.193000000000001
$ws.Range("C30").Value = -12.646
$ws.Range("B31").Value = 6.276999999999999
$ws.Range("A32").Value = -21.406
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -21.044
$ws.Range("D36").Value = -8.169
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.363
$ws.Range("C44").Value = -12.289
$ws.Range("B45").Value = 5.678000000000001
$ws.Range("A46").Value = -21.515
$ws.Range("A54").Value = -21.9
$ws.Range("A55").Value = -21.921
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.609
$ws.Range("D64").Value = -7.654999999999999
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.55
$ws.Range("C84").Value = -13.172
$ws.Range("C89").Value = -11.403
$ws.Range("D89").Value = -6.747999999999999
$ws.Range("A91").Value = -21.747
$ws.Range("C91").Value = -11.21
$ws.Range("C92").Value = -11.539
$ws.Range("A99").Value = -20.559
$ws.Range("B100").Value = 6.187
$ws.Range("B102").Value = 7.334999999999999
$ws.Range("C102").Value = -12.808
